$d = $word.ActiveDocument

# --- Step 1: remove the "Meta description: ..." paragraph that follows the
# title heading paragraph ("Play Arcane Gems for Free - Review").
$metaPara = $d.Paragraphs.Item(2)
if ($metaPara.Range.Text -like "Meta description*") {
  $metaPara.Range.Delete()
}

# --- Step 2: the final paragraph (the italic "Prompt: ..." image-prompt
# text) becomes the meta-description text, and a new bold heading-style
# paragraph ("Play Arcane Gems for Free - Review") is inserted right before
# it.

$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)

$oldPromptText = "Prompt: Create a cartoon-style feature image for the game " + [char]34 + "Arcane Gems" + [char]34 + " that features a happy Maya warrior with glasses. For this feature image, we want to bring in elements of both the theme of gems and the unique aspect of the game" + [char]39 + "s respin feature. The Maya warrior with glasses will add a touch of personality to the image and make it stand out. The Maya warrior should be depicted with a big smile on their face, eyes twinkling behind their glasses. They should be surrounded by piles of colorful gems, with one hand clutching a handful of gems, and the other hand pointing to the reels of the game. The reels should be shown on the image, with the game name " + [char]34 + "Arcane Gems" + [char]34 + " prominently displayed. The symbols on the reels should be vibrant and eye-catching, particularly highlighting the blue gem symbol that pays out the most. The symbol locking feature should also be represented, perhaps with the locked symbols being depicted as glowing and surrounded by a blue aura. Overall, the image should be lively and fun, drawing players in with its bright colors, charming character, and attention to the unique features of the game."
$newDescText = "Discover the pros and cons of Arcane Gems including its symbol locking and respin feature, high rewards, lack of wild symbol, and payout system. Play for free."

$find = $lastPara.Range.Find
$find.ClearFormatting()
$find.Text = $oldPromptText
$find.Replacement.ClearFormatting()
$find.Replacement.Text = $newDescText
[void]$find.Execute($null, $true, $false, $false, $false, $false, $true, 1, $false, $newDescText, 2)

# Insert a new paragraph right before this one for the bold title line.
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
[void]$lastPara.Range.InsertParagraphBefore()

$titlePara = $d.Paragraphs.Item($d.Paragraphs.Count - 1)
$r = $titlePara.Range
$r.Collapse(1)
[void]$r.InsertXML('<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Arcane Gems for Free - Review</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>')
